$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# 1) Insert a new "2022-Q1" sheet before "总计", cloned from "2021-Q3"
#    so it inherits the same (style-2) header/column formatting.
# ------------------------------------------------------------------
$src = $wb.Worksheets.Item(2)          # "2021-Q3"
$totalSheet = $wb.Worksheets.Item(3)   # "总计" (insertion point)
$src.Copy($totalSheet, $null)

$q1 = $wb.Worksheets.Item(3)
$q1.Name = "2022-Q1"

# Header row tweak specific to this sheet
$q1.Range("D1").Value = "基金规模"

# Row 2 - 005457 景顺长城量化小盘股票
$q1.Range("B2").Value = "'005457"
$q1.Range("C2").Value = "景顺长城量化小盘股票"
$q1.Range("D2").Value = "'9.49"
$q1.Range("E2").Value = "'93.39"
$q1.Range("F2").Value = "'1.99"
$q1.Range("G2").Value = "'0.1889"
$q1.Range("H2").Value = 2
$q1.Range("B2:G2").ClearFormats()

# Row 3 - 008072 景顺长城创业板综指增强
$q1.Range("A2").Copy()
$q1.Range("A3").PasteSpecial(-4122)
$q1.Range("A3").Value = 1
$q1.Range("B3").Value = "'008072"
$q1.Range("C3").Value = "景顺长城创业板综指增强"
$q1.Range("D3").Value = "'2.16"
$q1.Range("E3").Value = "'92.35"
$q1.Range("F3").Value = "'2.12"
$q1.Range("G3").Value = "'0.0458"
$q1.Range("H3").Value = 10
$q1.Range("B3:G3").ClearFormats()

# ------------------------------------------------------------------
# 2) "总计" sheet: prepend a 2022-Q1 summary row, pushing the
#    existing 2021-Q3 / 2021-Q2 rows down by one.
# ------------------------------------------------------------------
$tot = $wb.Worksheets.Item(4)

$oldB2 = $tot.Range("B2").Value2
$oldC2 = $tot.Range("C2").Value2
$oldD2 = $tot.Range("D2").Value2
$oldB3 = $tot.Range("B3").Value2
$oldC3 = $tot.Range("C3").Value2
$oldD3 = $tot.Range("D3").Value2

# Propagate the column-A styling to the (currently empty) row 4
$tot.Range("A3").Copy()
$tot.Range("A4").PasteSpecial(-4122)

$tot.Range("A4").Value = 2
$tot.Range("B4").Value = $oldB3
$tot.Range("C4").Value = $oldC3
$tot.Range("D4").Value = $oldD3

$tot.Range("A3").Value = 1
$tot.Range("B3").Value = $oldB2
$tot.Range("C3").Value = $oldC2
$tot.Range("D3").Value = $oldD2

$tot.Range("A2").Value = 0
$tot.Range("B2").Value = "2022-Q1"
$tot.Range("C2").Value = 2
$tot.Range("D2").Value = 0.23
